# Advent of Code 2025, day 4: append the new runtime sample to the
# RuntimesChart worksheet (row 6 -> Day 4, C# runtime 0.10504788s).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(6, 1).Value = 4
$ws.Cells.Item(6, 2).Value = 0.10504788

# Mirror Excel's normal "just typed a row" behaviour: selection moves to
# the newly entered row.
$ws.Range("A6:B6").Select()
